$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 38) below the existing table (which ends at row 37).
# Column A holds a date-like string ("2025/09/30") that must stay a plain text
# value (matching the existing rows), so force text format before assigning it
# and then drop back to the sheet's normal (unstyled) look, exactly like the
# other data rows.
$row = 38
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025/09/30"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = "火"
$ws.Range("C$row").Value = 6
$ws.Range("D$row").Value = 3
